# Fruta / hortaliza, semanal
# Insert a new weekly price record for Kiwi (Agricola del Norte S.A. de Arica)
# as a new row 20, pushing the existing rows 20-28 down to 21-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20; existing rows 20..28 shift down to 21..29
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new weekly record
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C20").Value = "Arica y Parinacota"
$ws.Range("D20").Value = 45134
$ws.Range("E20").Value = 15
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100101
$ws.Range("H20").Value = "Berries"
$ws.Range("I20").Value = 100101007
$ws.Range("J20").Value = "Kiwi"
$ws.Range("K20").Value = "Hayward"
$ws.Range("L20").Value = "Especial"
$ws.Range("M20").Value = 350
$ws.Range("N20").Value = 21000
$ws.Range("O20").Value = 22000
$ws.Range("P20").Value = 21429
$ws.Range("Q20").Value = "`$/bandeja 18 kilos"
$ws.Range("R20").Value = "Región de O'Higgins"
$ws.Range("S20").Value = 1190
$ws.Range("T20").Value = 18
